# Update Tab23 sheet: replace the stale/placeholder integer values in
# column N (row 13 header row + data rows 23..97) with the corrected
# decimal figures, and shrink the saved window height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab23")

# Map of cell reference -> corrected value.
$values = [ordered]@{
    "N13" = 2525.2400042378199
    "N23" = 2373.8133079182198
    "N38" = 2258.1855425952399
    "N45" = 3336.2071632669099
    "N61" = 2590.5500867576102
    "N62" = 2597.67579128335
    "N63" = 2973.5091189032701
    "N64" = 3047.2970802363898
    "N65" = 2804.1083144570098
    "N66" = 2916.5117556755899
    "N67" = 2511.50676411314
    "N68" = 2755.4966645743498
    "N69" = 2230.1906916993898
    "N70" = 2320.2769207746501
    "N71" = 2590.5500867576102
    "N72" = 2262.22687253816
    "N73" = 2437.6284547318101
    "N74" = 3352.0306002006701
    "N75" = 2832.1198234047301
    "N76" = 3051.80122409674
    "N77" = 3450.4861828215398
    "N78" = 2565.1844073995699
    "N79" = 3002.9781012830599
    "N80" = 2612.2435405891802
    "N81" = 2971.0449106783599
    "N82" = 2357.5412526099699
    "N83" = 2048.36264164764
    "N84" = 2704.5365449258402
    "N85" = 2560.2928867453902
    "N86" = 2938.7327190576202
    "N87" = 3136.9501503757401
    "N88" = 3430.63190255068
    "N89" = 2365.7601031254299
    "N90" = 2515.0144422708399
    "N91" = 2586.6169510907398
    "N92" = 2782.3243861455298
    "N93" = 2343.3037766356401
    "N94" = 2683.4725649816901
    "N95" = 2383.6527428618601
    "N96" = 2536.0402697965901
    "N97" = 2235.95057650024
}

foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# Shrink the saved window height (xWindow/yWindow/windowWidth untouched).
$win = $wb.Windows.Item(1)
$win.Height = 11490
